$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "43.213.76"
    "E2" = "  +0.70%  "
    "D3" = "2.306.52"
    "E3" = "  +0.53%  "
    "E4" = "  +0.01%  "
    "D5" = "300.95"
    "E5" = "  +0.08%  "
    "D6" = "97.91"
    "E6" = "  -1.08%  "
    "E7" = "  +2.90%  "
    "E8" = "  -0.03%  "
    "D9" = "0.517"
    "E9" = "  +1.12%  "
    "D10" = "36.32"
    "E10" = "  +0.74%  "
    "E11" = "  +0.38%  "
    "E12" = "  +0.74%  "
    "D13" = "17.69"
    "E13" = "  -3.79%  "
    "D14" = "6.90"
    "E14" = "  -0.81%  "
    "D15" = "2.665.17"
    "E15" = "  +0.55%  "
    "D16" = "2.348.40"
    "E16" = "  +2.94%  "
    "E17" = "  -1.20%  "
    "D18" = "43.076.82"
    "E18" = "  +0.61%  "
    "D19" = "13.08"
    "E19" = "  +4.60%  "
    "D20" = "0.0₃0912"
    "E20" = "  +1.01%  "
    "E21" = "  +0.50%  "
    "D22" = "68.32"
    "E22" = "  +0.94%  "
    "D23" = "238.20"
    "E23" = "  +1.22%  "
    "D24" = "2.22"
    "E24" = "  -0.38%  "
    "E25" = "  -0.52%  "
    "D26" = "2.42"
    "E26" = "  -0.70%  "
    "E27" = "  +0.00%  "
    "D28" = "25.36"
    "E28" = "  +1.46%  "
    "D29" = "9.18"
    "E29" = "  +0.70%  "
    "E30" = "  -13.52%  "
    "D31" = "162.97"
    "E31" = "  -2.44%  "
    "D32" = "33.29"
    "E32" = "  -3.33%  "
    "D33" = "0.999"
    "E33" = "  -0.02%  "
    "E34" = "  +2.91%  "
    "D35" = "18.24"
    "E35" = "  +3.18%  "
    "D36" = "4.75"
    "E36" = "  +1.83%  "
    "D38" = "0.0695"
    "E38" = "  +1.16%  "
    "E39" = "  +1.53%  "
    "E40" = "  +0.16%  "
    "E42" = "  -1.73%  "
    "D43" = "2.016.10"
    "E43" = "  +1.98%  "
    "E44" = "  -0.74%  "
    "D45" = "2.19"
    "E45" = "  -6.96%  "
    "D46" = "10.25"
    "E46" = "  +1.60%  "
    "E47" = "  +0.92%  "
    "D48" = "2.86"
    "E48" = "  -0.83%  "
    "D49" = "54.44"
    "E49" = "  -1.23%  "
    "D50" = "2.536.81"
    "E50" = "  +0.78%  "
    "E51" = "  +0.26%  "
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    $range.ClearFormats()
}
